$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Date of Birth cell as Text so the value "1232-10-23" is
# stored as a literal string instead of being auto-converted into a date
# serial number.
$ws.Range("C5").NumberFormat = "@"

$ws.Range("A5").Value = "P1004"
$ws.Range("B5").Value = "Hi"
$ws.Range("C5").Value = "1232-10-23"
$ws.Range("D5").Value = "Female"
$ws.Range("E5").Value = "O-"
$ws.Range("F5").Value = 98736151
$ws.Range("G5").Value = "1a1"
